$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Title "Documentation accéssibilité sur téléphone" gets split into
# three runs, with "accéssibilité" wrapped in spellcheck proofErr markers
# (Word flags it as a misspelling), while keeping identical run formatting.
# ---------------------------------------------------------------------------
$titleRng = $d.Content
$titleFound = $titleRng.Find.Execute("Documentation accéssibilité sur téléphone", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($titleFound) {
    $titleFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="259" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI Light" w:hAnsi="Segoe UI Light" w:cs="Segoe UI Light"/><w:color w:val="009FE3"/><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Light" w:hAnsi="Segoe UI Light" w:cs="Segoe UI Light"/><w:color w:val="009FE3"/><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr><w:t xml:space="preserve">Documentation </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Light" w:hAnsi="Segoe UI Light" w:cs="Segoe UI Light"/><w:color w:val="009FE3"/><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr><w:t>accéssibilité</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Light" w:hAnsi="Segoe UI Light" w:cs="Segoe UI Light"/><w:color w:val="009FE3"/><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr><w:t xml:space="preserve"> sur téléphone</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $titleRng.InsertXML($titleFrag) | Out-Null
}

# ---------------------------------------------------------------------------
# Change 2: Second page heading ".Page d'accueil" (after the "2" run) gets
# renamed to ".Page " + "Evènements" (two runs, same formatting as before),
# and a new blank paragraph (matching the existing trailing blank paragraph's
# properties) is inserted right after it.
# ---------------------------------------------------------------------------
$pageRng = $d.Content
$pageFound = $pageRng.Find.Execute("2.Page d’accueil", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($pageFound) {
    $pageFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="259" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Segoe UI Black" w:hAnsi="Segoe UI Black" w:cs="Segoe UI Semibold"/><w:color w:val="AFCA0B"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Black" w:hAnsi="Segoe UI Black" w:cs="Segoe UI Semibold"/><w:color w:val="AFCA0B"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Black" w:hAnsi="Segoe UI Black" w:cs="Segoe UI Semibold"/><w:color w:val="AFCA0B"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">.Page </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI Black" w:hAnsi="Segoe UI Black" w:cs="Segoe UI Semibold"/><w:color w:val="AFCA0B"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>Evènements</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="259" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Segoe UI Black" w:hAnsi="Segoe UI Black" w:cs="Segoe UI Semibold"/><w:color w:val="AFCA0B"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    $pageRng.InsertXML($pageFrag) | Out-Null
}
